$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Range("B6").Value = "US Dollar"
$ws.Range("A6").Value = "currency"
$ws.Range("B6").Interior.Color = $ws.Range("B1").Interior.Color
$ws.Range("B6").Font.Name = $ws.Range("B1").Font.Name
$ws.Range("B6").Font.Size = $ws.Range("B1").Font.Size
$ws.Range("B6").Font.Bold = $ws.Range("B1").Font.Bold
$ws.Columns.Item(3).Delete()
